# Weekly update: a new price-report row for
# "Terminal Hortofrutícola Agro Chillán - Mango" is inserted as the new
# most-recent record (row 9), pushing all existing data rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (row 1 is the header, row 2-48 hold prior
# weekly records sorted most-recent-first) - this shifts old rows 9..48
# down to 10..49 and extends the sheet's used range to T49.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row with this week's reading.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44453
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100108
$ws.Range("H9").Value = "Tropicales y subtropicales"
$ws.Range("I9").Value = 100108002
$ws.Range("J9").Value = "Mango"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 8500
$ws.Range("O9").Value = 9000
$ws.Range("P9").Value = 8750
$ws.Range("Q9").Value = "$/bandeja 4 kilos"
$ws.Range("R9").Value = "Brasil"
$ws.Range("S9").Value = 2188
$ws.Range("T9").Value = 4
